$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The mods:note tag in V2 dropped its type="local" attribute -- it now reads
# <mods:note displayLabel="Description"> instead of
# <mods:note displayLabel="Description" type="local">
$ws.Range("V2").Value = '<mods:note displayLabel="Description">'

[void]$ws.Range("V2").Select()
